$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new blank rows before row 12. This pushes the existing
# "Programa resumido:" ... "Requisitos:"/"LOM3037..." block (old rows 12-22)
# down to rows 15-25, matching the target layout.
$ws.Rows("12:14").Insert()

# Row-insert in this runtime copies column-A formatting down from the row
# above into the new rows' column A (so A12/A13/A14 all pick up bold style).
# Only A12 should end up populated/formatted; clear the spurious A13/A14
# cells entirely so they don't linger as empty-but-styled cells.
$ws.Range("A13:A14").Clear()

# New row 12: section header "Docentes responsáveis:" in column A, styled
# like the other section headers (copy format from the header that is about
# to land on row 15, "Programa resumido:").
$ws.Range("A15").Copy()
$ws.Range("A12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"

# New rows 13-14: the two instructors, duplicated into column B (current)
# and column C (modified/highlighted in red), matching the existing
# current-vs-modified layout used throughout the sheet. Copy formatting from
# the B/C cells landing on row 15.
$ws.Range("B15").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("C15").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)

$ws.Cells.Item(13, 2).Value = "3577649 - Carlos Angelo Nunes"
$ws.Cells.Item(13, 3).Value = "3577649 - Carlos Angelo Nunes"

$ws.Cells.Item(14, 2).Value = "1922320 - Sebastiao Ribeiro"
$ws.Cells.Item(14, 3).Value = "1922320 - Sebastiao Ribeiro"
